# Update market-price derived columns (H-N) across leve-profit sheets
# per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3832.6667
$ws.Range("H72").Value = 3832.6667
$ws.Range("H76").Value = 3074.558
$ws.Range("I76").Value = 3000.08
$ws.Range("K76").Value = 3000.08
$ws.Range("M76").Value = -2685.08
$ws.Range("H79").Value = 3074.558
$ws.Range("I79").Value = 3000.08
$ws.Range("K79").Value = 3000.08
$ws.Range("M79").Value = -1908.08
$ws.Range("H80").Value = 11448339
$ws.Range("I80").Value = 19608432
$ws.Range("J80").Value = 777448.4
$ws.Range("K80").Value = 58825296
$ws.Range("L80").Value = 2332345.2
$ws.Range("M80").Value = -58824298
$ws.Range("N80").Value = -2334341.2
$ws.Range("H82").Value = 66669136
$ws.Range("I82").Value = 2086.0833
$ws.Range("J82").Value = 333337340
$ws.Range("K82").Value = 6258.249899999999
$ws.Range("L82").Value = 1000012020
$ws.Range("M82").Value = -5852.249899999999
$ws.Range("N82").Value = -1000012832
$ws.Range("H83").Value = 11448339
$ws.Range("I83").Value = 19608432
$ws.Range("J83").Value = 777448.4
$ws.Range("K83").Value = 176475888
$ws.Range("L83").Value = 6997035.600000001
$ws.Range("M83").Value = -176470896
$ws.Range("N83").Value = -7007019.600000001
$ws.Range("H85").Value = 66669136
$ws.Range("I85").Value = 2086.0833
$ws.Range("J85").Value = 333337340
$ws.Range("K85").Value = 6258.249899999999
$ws.Range("L85").Value = 1000012020
$ws.Range("M85").Value = -4854.249899999999
$ws.Range("N85").Value = -1000014828
$ws.Range("H88").Value = 6768.8
$ws.Range("I88").Value = 9800
$ws.Range("J88").Value = 3737.6
$ws.Range("K88").Value = 9800
$ws.Range("L88").Value = 3737.6
$ws.Range("M88").Value = -9394
$ws.Range("N88").Value = -4549.6
$ws.Range("H91").Value = 6768.8
$ws.Range("I91").Value = 9800
$ws.Range("J91").Value = 3737.6
$ws.Range("K91").Value = 9800
$ws.Range("L91").Value = 3737.6
$ws.Range("M91").Value = -8396
$ws.Range("N91").Value = -6545.6
$ws.Range("H132").Value = 22103.088
$ws.Range("I132").Value = 3458.0264
$ws.Range("J132").Value = 110667.125
$ws.Range("K132").Value = 10374.0792
$ws.Range("L132").Value = 332001.375
$ws.Range("M132").Value = -7844.0792
$ws.Range("N132").Value = -337061.375
$ws.Range("H138").Value = 1547.98
$ws.Range("I138").Value = 770.34283
$ws.Range("J138").Value = 1966.7076
$ws.Range("K138").Value = 2311.02849
$ws.Range("L138").Value = 5900.1228
$ws.Range("M138").Value = 2828.97151
$ws.Range("N138").Value = -16180.1228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7726.9165
$ws.Range("I32").Value = 7007.56
$ws.Range("K32").Value = 7007.56
$ws.Range("M32").Value = -6720.56
$ws.Range("H74").Value = 1458.3726
$ws.Range("I74").Value = 1215.4546
$ws.Range("J74").Value = 2985.2856
$ws.Range("K74").Value = 1215.4546
$ws.Range("L74").Value = 2985.2856
$ws.Range("M74").Value = -341.4546
$ws.Range("N74").Value = -4733.2856
$ws.Range("H77").Value = 1458.3726
$ws.Range("I77").Value = 1215.4546
$ws.Range("J77").Value = 2985.2856
$ws.Range("K77").Value = 6077.273
$ws.Range("L77").Value = 14926.428
$ws.Range("M77").Value = -1709.273
$ws.Range("N77").Value = -23662.428
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1930.6666
$ws.Range("I86").Value = 1896.5
$ws.Range("J86").Value = 1999
$ws.Range("K86").Value = 1896.5
$ws.Range("L86").Value = 1999
$ws.Range("M86").Value = -773.5
$ws.Range("N86").Value = -4245
$ws.Range("H89").Value = 1930.6666
$ws.Range("I89").Value = 1896.5
$ws.Range("J89").Value = 1999
$ws.Range("K89").Value = 9482.5
$ws.Range("L89").Value = 9995
$ws.Range("M89").Value = -3866.5
$ws.Range("N89").Value = -21227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4322.9443
$ws.Range("I31").Value = 2034.2354
$ws.Range("K31").Value = 2034.2354
$ws.Range("M31").Value = -1739.2354
$ws.Range("H34").Value = 4322.9443
$ws.Range("I34").Value = 2034.2354
$ws.Range("K34").Value = 2034.2354
$ws.Range("M34").Value = -1832.2354
$ws.Range("H132").Value = 57694.32
$ws.Range("I132").Value = 953.5
$ws.Range("J132").Value = 203599.28
$ws.Range("K132").Value = 2860.5
$ws.Range("L132").Value = 610797.84
$ws.Range("M132").Value = -330.5
$ws.Range("N132").Value = -615857.84
$ws.Range("H141").Value = 20232.223
$ws.Range("J141").Value = 6636.25
$ws.Range("L141").Value = 6636.25
$ws.Range("N141").Value = -16996.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 19068040
$ws.Range("I33").Value = 199.33333
$ws.Range("J33").Value = 24268360
$ws.Range("K33").Value = 1195.99998
$ws.Range("L33").Value = 145610160
$ws.Range("M33").Value = -912.9999800000001
$ws.Range("N33").Value = -145610726
$ws.Range("H113").Value = 3048.9048
$ws.Range("I113").Value = 4270.815
$ws.Range("J113").Value = 849.4666999999999
$ws.Range("K113").Value = 12812.445
$ws.Range("L113").Value = 2548.4001
$ws.Range("M113").Value = -10642.445
$ws.Range("N113").Value = -6888.4001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3234.889
$ws.Range("I102").Value = 4930.6665
$ws.Range("J102").Value = 2387
$ws.Range("K102").Value = 4930.6665
$ws.Range("L102").Value = 2387
$ws.Range("M102").Value = -3308.6665
$ws.Range("N102").Value = -5631

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1227.8572
$ws.Range("I55").Value = 1365
$ws.Range("J55").Value = 1125
$ws.Range("K55").Value = 1365
$ws.Range("L55").Value = 1125
$ws.Range("M55").Value = -1192
$ws.Range("N55").Value = -1471
$ws.Range("H82").Value = 1689.091
$ws.Range("I82").Value = 1730
$ws.Range("J82").Value = 1640
$ws.Range("K82").Value = 1730
$ws.Range("L82").Value = 1640
$ws.Range("M82").Value = -1369
$ws.Range("N82").Value = -2362
$ws.Range("H85").Value = 1689.091
$ws.Range("I85").Value = 1730
$ws.Range("J85").Value = 1640
$ws.Range("K85").Value = 1730
$ws.Range("L85").Value = 1640
$ws.Range("M85").Value = -482
$ws.Range("N85").Value = -4136
$ws.Range("H122").Value = 102268
$ws.Range("I122").Value = 145041.72
$ws.Range("J122").Value = 2462.6667
$ws.Range("K122").Value = 435125.16
$ws.Range("L122").Value = 7388.000100000001
$ws.Range("M122").Value = -432675.16
$ws.Range("N122").Value = -12288.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1731743.6
$ws.Range("I126").Value = 1839883.9
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 5519651.699999999
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -5517181.699999999
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 1349.4762
$ws.Range("I132").Value = 1168.8334
$ws.Range("J132").Value = 2433.3333
$ws.Range("K132").Value = 3506.5002
$ws.Range("L132").Value = 7299.999899999999
$ws.Range("M132").Value = -976.5001999999999
$ws.Range("N132").Value = -12359.9999
$ws.Range("H141").Value = 38545.273
$ws.Range("J141").Value = 38545.273
$ws.Range("L141").Value = 38545.273
$ws.Range("N141").Value = -48905.273
